$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected/active cell (also resets view to no longer be scrolled to G1)
$ws.Range("F5").Select()

# Row 4 updates
$ws.Range("C4").Value = "cyclone"
$ws.Range("D4").Value = "luminaire"
$ws.Range("E4").Value = "luminaire"
$ws.Range("K4").Value = 30
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 30

# Row 5 addition
$ws.Range("C5").Value = "wind slash"

# Row 6 addition
$ws.Range("C6").Value = "lightning 2"
